$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.567.43'
$ws.Range('E2').Value = '  +2.36%  '
$ws.Range('D3').Value = '3.074.39'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').Value = "'516.44"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.66%  '
$ws.Range('D6').Value = "'141.46"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.31%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').Value = "'0.433"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  +2.78%  '
$ws.Range('E10').Value = '  +1.25%  '
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').Value = '3.596.51'
$ws.Range('E12').Value = '  +1.35%  '
$ws.Range('D13').Value = "'0.128"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.59%  '
$ws.Range('D14').Value = "'25.61"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('E15').Value = '  +0.88%  '
$ws.Range('D16').Value = '57.582.93'
$ws.Range('E16').Value = '  +2.24%  '
$ws.Range('D17').Value = '3.076.71'
$ws.Range('E17').Value = '  +1.53%  '
$ws.Range('D18').Value = "'6.09"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('D19').Value = "'13.01"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('E20').Value = '  +2.30%  '
$ws.Range('D21').Value = "'334.90"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.32%  '
$ws.Range('D22').Value = "'0.999"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('E23').Value = '  +1.20%  '
$ws.Range('D24').Value = "'65.82"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.17%  '
$ws.Range('E25').Value = '  +4.18%  '
$ws.Range('D26').Value = "'1.01"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').Value = '0.0₃0915'
$ws.Range('E27').Value = '  +3.39%  '
$ws.Range('D28').Value = "'6.36"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.20%  '
$ws.Range('D29').Value = "'7.18"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.46%  '
$ws.Range('D30').Value = "'1.82"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.99%  '
$ws.Range('E31').Value = '  +2.62%  '
$ws.Range('D32').Value = "'1.17"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('D33').Value = "'154.26"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = "'4.48"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.48%  '
$ws.Range('B35').Value = 'EnergySwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D35').Value = "'27.13"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.21%  '
$ws.Range('E36').Value = '  +1.26%  '
$ws.Range('D37').Value = "'1.28"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.31%  '
$ws.Range('D38').Value = "'0.0676"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.33%  '
$ws.Range('D39').Value = '3.109.41'
$ws.Range('E39').Value = '  +1.74%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = "'3.91"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.42%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = "'36.99"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.06%  '
$ws.Range('D42').Value = "'0.999"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('D43').Value = "'0.656"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').Value = '2.264.54'
$ws.Range('E44').Value = '  +4.02%  '
$ws.Range('D45').Value = "'0.0256"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.71%  '
$ws.Range('D46').Value = "'1.38"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.11%  '
$ws.Range('D47').Value = "'19.98"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('D48').Value = "'0.928"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('D49').Value = "'5.86"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('D50').Value = "'263.59"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +15.63%  '
$ws.Range('E51').Value = '  +1.97%  '
